# Natmi following Dr Hou advice
# Re-derive LR-pair table adding "ECs" as a third sending/target cluster
# alongside the existing "FAPs" and "sCs", expanding the data grid from
# 2 sending clusters x 3 target clusters (6 rows) to 3 x 3 (9 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @("ECs", "Il4", "Il4ra", "ECs", 2, 0.6666666666666666, 0.6580480000000001, 1.974144, 0.1830904640197835, 0.1830904640197835, 3, 1, 49.00494433333333, 147.014833, 0.6271308668609277, 0.6271308668609277, 32.24760560866134, 290.228450477952, 0.1148216814146963, 0.1148216814146963),
    @("ECs", "Il4", "Il4ra", "FAPs", 2, 0.6666666666666666, 0.6580480000000001, 1.974144, 0.1830904640197835, 0.1830904640197835, 3, 1, 26.17647666666667, 78.52943, 0.3349881675544572, 0.3349881675544572, 17.22537811754667, 155.02840305792, 0.06133313903868254, 0.06133313903868254),
    @("ECs", "Il4", "Il4ra", "sCs", 2, 0.6666666666666666, 0.6580480000000001, 1.974144, 0.1830904640197835, 0.1830904640197835, 3, 1, 2.960075333333334, 8.880226, 0.03788096558461519, 0.03788096558461519, 1.947871652949334, 17.530844876544, 0.006935643566404642, 0.006935643566404642),
    @("FAPs", "Il4", "Il4ra", "ECs", 3, 1, 1.115367, 3.346101, 0.3103315587652478, 0.3103315587652479, 3, 1, 49.00494433333333, 147.014833, 0.6271308668609277, 0.6271308668609277, 54.658497746237, 491.926479716133, 0.1946184994627528, 0.1946184994627528),
    @("FAPs", "Il4", "Il4ra", "FAPs", 3, 1, 1.115367, 3.346101, 0.3103315587652478, 0.3103315587652479, 3, 1, 26.17647666666667, 78.52943, 0.3349881675544572, 0.3349881675544572, 29.19637825027, 262.76740425243, 0.1039574002050887, 0.1039574002050887),
    @("FAPs", "Il4", "Il4ra", "sCs", 3, 1, 1.115367, 3.346101, 0.3103315587652478, 0.3103315587652479, 3, 1, 2.960075333333334, 8.880226, 0.03788096558461519, 0.03788096558461519, 3.301570344314, 29.714133098826, 0.01175565909740634, 0.01175565909740634),
    @("sCs", "Il4", "Il4ra", "ECs", 3, 1, 1.820699, 5.462097, 0.5065779772149687, 0.5065779772149687, 3, 1, 49.00494433333333, 147.014833, 0.6271308668609277, 0.6271308668609277, 89.22325314275567, 803.0092782848011, 0.3176906859834785, 0.3176906859834785),
    @("sCs", "Il4", "Il4ra", "FAPs", 3, 1, 1.820699, 5.462097, 0.5065779772149687, 0.5065779772149687, 3, 1, 26.17647666666667, 78.52943, 0.3349881675544572, 0.3349881675544572, 47.65948489052334, 428.93536401471, 0.1696976283106859, 0.1696976283106859),
    @("sCs", "Il4", "Il4ra", "sCs", 3, 1, 1.820699, 5.462097, 0.5065779772149687, 0.5065779772149687, 3, 1, 2.960075333333334, 8.880226, 0.03788096558461519, 0.03788096558461519, 5.389406199324667, 48.50465579392201, 0.0191896629208042, 0.0191896629208042)
)

$startRow = 2
for ($i = 0; $i -lt $newData.Length; $i++) {
    $rowVals = $newData[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowVals[$j]
    }
}

Write-Output "Updated $($newData.Length) data rows (rows $startRow..$($startRow + $newData.Length - 1))."
